$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price (D) and Volume(1h) (E) columns so that
# numeric-looking strings keep their original text formatting (trailing
# zeros, multi-dot thousand separators, etc.) instead of being coerced
# into Excel numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.910.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.747.14'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.04'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.28%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.98%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.382'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.50'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -18.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.231.55'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.50'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.592.12'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000149'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.751.43'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.13'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '354.57'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.75'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.537'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.17'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.40'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0892'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.37%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.95'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.57'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.20'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.98%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.28%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.15'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +7.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.13'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '324.11'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.92'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.70%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0587'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.21'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0253'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '134.43'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.622'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.100'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.58%  '
